$wb = $excel.ActiveWorkbook

# --- Add the new "CarDetails" worksheet after the last existing sheet ---
$sheetCount = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($sheetCount)
$bikeSheet  = $wb.Worksheets.Item(1)

$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "CarDetails"

# --- Table data (header + 13 upcoming Tata car rows) ---
$data = @(
    @("Car Name",              "Car Price",       "Launch Date"),
    @("Tata Altroz Racer",     "Rs. 10.00 Lakh",  "20 Mar 2024"),
    @("Tata Curvv EV",         "Rs. 20.00 Lakh",  "Jul 2024"),
    @("Tata Curvv",            "Rs. 10.50 Lakh",  "Aug 2024"),
    @("Tata Avinya",           "Rs. 30.00 Lakh",  "Jan 2025"),
    @("Tata Harrier EV",       "Rs. 30.00 Lakh",  "Apr 2025"),
    @("Tata Punch 2025",       "Rs. 6.00 Lakh",   "Jun 2025"),
    @("Tata Sierra",           "Rs. 25.00 Lakh",  "Dec 2025"),
    @("Tata Kite 5",           "Rs. 4.50 Lakh",   "Unrevealed"),
    @("Tata Atmos",            "Rs. 12.00 Lakh",  "Unrevealed"),
    @("Tata H7X",              "Rs. 15.00 Lakh",  "Unrevealed"),
    @("Tata Altroz EV",        "Rs. 14.00 Lakh",  "Unrevealed"),
    @("Tata Hexa",             "Rs. 14.00 Lakh",  "Unrevealed"),
    @("Tata EVision Electric", "Rs. 12.00 Lakh",  "Unrevealed")
)

# A handful of the "Month Year" launch-date strings (e.g. "Jul 2024") look
# like dates to Excel's smart-entry and would otherwise get silently turned
# into date serial numbers. Mark those specific cells as Text first so the
# literal string is preserved, then strip the temporary format again so the
# cell ends up back on the sheet's normal (default) style.
$dateRiskValues = @("Jul 2024", "Aug 2024", "Jan 2025", "Apr 2025", "Jun 2025", "Dec 2025")

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 1
    $ws.Cells.Item($rowNum, 1).Value = $data[$r][0]
    $ws.Cells.Item($rowNum, 2).Value = $data[$r][1]

    $cellC = $ws.Cells.Item($rowNum, 3)
    $value = $data[$r][2]
    if ($dateRiskValues -contains $value) {
        $cellC.NumberFormat = "@"
        $cellC.Value = $value
        # Reset back to a plain/default-styled cell (no lingering text format)
        $bikeSheet.Range("A4").Copy()
        $cellC.PasteSpecial(-4122)
    } else {
        $cellC.Value = $value
    }
}

# --- Header styling: reuse the same "white text on coloured fill" look
# that the BikeDetails header row already uses ---
$bikeSheet.Range("A1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$bikeSheet.Range("B1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$bikeSheet.Range("C1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# --- Column widths roughly matching the content (bestFit-style) ---
$ws.Columns.Item(1).ColumnWidth = 17.04
$ws.Columns.Item(2).ColumnWidth = 12.37
$ws.Columns.Item(3).ColumnWidth = 11.14

# Restore the original active sheet/selection (BikeDetails stays active)
$bikeSheet.Activate()
